# Update the dSF (column F) values on Sheet1 to reflect the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    3  = -9
    4  = -2
    5  = -6
    6  = -5
    7  = -2
    9  = -1
    12 = 2
    14 = -5
    16 = 1
    18 = 11
    19 = -2
    21 = -1
    24 = -8
    27 = -8
    28 = -8
    29 = -10
    31 = 7
    34 = -8
    35 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
